# Rewrite the single-sheet test-data workbook from the "CURA Healthcare" booking
# fixture (TCID/Username/Password/Facility/Visit Date/Comment/Title, 3 rows) to the
# smaller Login fixture (TCID/Username/Password, 2 data rows) used by the Login tests,
# and add a thin border around the new table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the columns/rows that are no longer part of the table -----------------
# Clear content + formatting of everything outside the new A:C table first (this
# drops the stray styled cells in columns D-G and in row 4 before we physically
# remove the now-empty columns/row).
$ws.Range("D1:G4").Clear()

# Row 4 (old TC3 / "CURA Healthcare" row) is removed entirely.
$ws.Rows.Item(4).Delete()

# Columns D:F (Facility / Visit Date / Comment) and the old column G (Title) are
# removed. Column K's leftover width definition shifts down to column G, matching
# the target layout.
$ws.Range("D1:F1").EntireColumn.Delete()
$ws.Range("G1").EntireColumn.Delete()

# --- New table content --------------------------------------------------------------
# Header row (keeps the existing bold / filled header style already present on row 1).
$ws.Range("A1").Value = "TCID"
$ws.Range("A2").Value = "TC1"
$ws.Range("A3").Value = "TC2"
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("C3").Value = "invalidPassword"
$ws.Range("B2").Value = "johndoe"
$ws.Range("B3").Value = "johndoe"
$ws.Range("C2").Value = "validPassword1!"

# --- Formatting: add a thin border around the header and the data rows --------------
$ws.Range("A1:C1").Borders.LineStyle = 1
$ws.Range("A2:C3").Borders.LineStyle = 1

# --- Selection -----------------------------------------------------------------------
[void]$ws.Range("C1").Select()

Write-Host "Rebuilt Sheet1 login test-data table (TCID/Username/Password)."
